$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6 and add new rows 7-11 with refreshed TPM-derived values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Oprm1"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.8377936666666667
$ws.Range("N2").Value = 2.513381
$ws.Range("O2").Value = 0.7130909380817101
$ws.Range("P2").Value = 0.7130909380817101
$ws.Range("Q2").Value = 121.2985359680683
$ws.Range("R2").Value = 1091.686823712615
$ws.Range("S2").Value = 0.1732844020743374
$ws.Range("T2").Value = 0.1789948690310812

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Oprm1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 144.783305
$ws.Range("H3").Value = 434.349915
$ws.Range("I3").Value = 0.2430046335191003
$ws.Range("J3").Value = 0.251012682214973
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3370826666666667
$ws.Range("N3").Value = 1.011248
$ws.Range("O3").Value = 0.2869090619182899
$ws.Range("P3").Value = 0.2869090619182899
$ws.Range("Q3").Value = 48.80394253821333
$ws.Range("R3").Value = 439.23548284392
$ws.Range("S3").Value = 0.0697202314447629
$ws.Range("T3").Value = 0.07201781318389169

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Oprm1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 82.24887099999999
$ws.Range("H4").Value = 246.746613
$ws.Range("I4").Value = 0.1380466950572427
$ws.Range("J4").Value = 0.1425959278859072
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.8377936666666667
$ws.Range("N4").Value = 2.513381
$ws.Range("O4").Value = 0.7130909380817101
$ws.Range("P4").Value = 0.7130909380817101
$ws.Range("Q4").Value = 68.90758321428366
$ws.Range("R4").Value = 620.168248928553
$ws.Range("S4").Value = 0.09843984727744895
$ws.Range("T4").Value = 0.1016838639827934

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Oprm1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 82.24887099999999
$ws.Range("H5").Value = 246.746613
$ws.Range("I5").Value = 0.1380466950572427
$ws.Range("J5").Value = 0.1425959278859072
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3370826666666667
$ws.Range("N5").Value = 1.011248
$ws.Range("O5").Value = 0.2869090619182899
$ws.Range("P5").Value = 0.2869090619182899
$ws.Range("Q5").Value = 27.72466876700267
$ws.Range("R5").Value = 249.522018903024
$ws.Range("S5").Value = 0.03960684777979372
$ws.Range("T5").Value = 0.04091206390311374

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Oprm1"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 163.8590903333333
$ws.Range("H6").Value = 491.577271
$ws.Range("I6").Value = 0.2750214756820535
$ws.Range("J6").Value = 0.284084617144743
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.8377936666666667
$ws.Range("N6").Value = 2.513381
$ws.Range("O6").Value = 0.7130909380817101
$ws.Range("P6").Value = 0.7130909380817101
$ws.Range("Q6").Value = 137.2801081070279
$ws.Range("R6").Value = 1235.520972963251
$ws.Range("S6").Value = 0.1961153220867317
$ws.Range("T6").Value = 0.2025781661343282

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Oprm1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 163.8590903333333
$ws.Range("H7").Value = 491.577271
$ws.Range("I7").Value = 0.2750214756820535
$ws.Range("J7").Value = 0.284084617144743
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3370826666666667
$ws.Range("N7").Value = 1.011248
$ws.Range("O7").Value = 0.2869090619182899
$ws.Range("P7").Value = 0.2869090619182899
$ws.Range("Q7").Value = 55.23405912713423
$ws.Range("R7").Value = 497.1065321442081
$ws.Range("S7").Value = 0.07890615359532173
$ws.Range("T7").Value = 0.08150645101041473

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Oprm1"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 57.0238095
$ws.Range("H8").Value = 114.047619
$ws.Range("I8").Value = 0.09570889357312636
$ws.Range("J8").Value = 0.06590860906562239
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.8377936666666667
$ws.Range("N8").Value = 2.513381
$ws.Range("O8").Value = 0.7130909380817101
$ws.Range("P8").Value = 0.7130909380817101
$ws.Range("Q8").Value = 47.7741864483065
$ws.Range("R8").Value = 286.645118689839
$ws.Range("S8").Value = 0.06824914470082323
$ws.Range("T8").Value = 0.04699883186626537

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Oprm1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 57.0238095
$ws.Range("H9").Value = 114.047619
$ws.Range("I9").Value = 0.09570889357312636
$ws.Range("J9").Value = 0.06590860906562239
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3370826666666667
$ws.Range("N9").Value = 1.011248
$ws.Range("O9").Value = 0.2869090619182899
$ws.Range("P9").Value = 0.2869090619182899
$ws.Range("Q9").Value = 19.221737769752
$ws.Range("R9").Value = 115.330426618512
$ws.Range("S9").Value = 0.02745974887230312
$ws.Range("T9").Value = 0.01890977719935701

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Oprm1"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 147.8896333333333
$ws.Range("H10").Value = 443.6689
$ws.Range("I10").Value = 0.2482183021684772
$ws.Range("J10").Value = 0.2563981636887546
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.8377936666666667
$ws.Range("N10").Value = 2.513381
$ws.Range("O10").Value = 0.7130909380817101
$ws.Range("P10").Value = 0.7130909380817101
$ws.Range("Q10").Value = 123.9009981723222
$ws.Range("R10").Value = 1115.1089835509
$ws.Range("S10").Value = 0.1770022219423688
$ws.Range("T10").Value = 0.1828352070672419

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Oprm1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 147.8896333333333
$ws.Range("H11").Value = 443.6689
$ws.Range("I11").Value = 0.2482183021684772
$ws.Range("J11").Value = 0.2563981636887546
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.3370826666666667
$ws.Range("N11").Value = 1.011248
$ws.Range("O11").Value = 0.2869090619182899
$ws.Range("P11").Value = 0.2869090619182899
$ws.Range("Q11").Value = 49.85103197635556
$ws.Range("R11").Value = 448.6592877872001
$ws.Range("S11").Value = 0.07121608022610841
$ws.Range("T11").Value = 0.0735629566215127

